$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '326.63'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-1.28%'
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '44.30'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-0.26%'
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.289'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-4.49%'
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08350'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '2.30%'
# Row 6
$ws.Range("B6").Value = 'FTXToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.938'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-6.14%'
# Row 7
$ws.Range("B7").Value = 'MXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9685'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-0.60%'
# Row 8
$ws.Range("B8").Value = 'BTSEToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.541'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-3.92%'
# Row 9
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1127'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '1.02%'
# Row 10
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1893'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '0.08%'
# Row 11
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09681'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-3.62%'
# Row 12
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.04595'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-3.31%'
# Row 13
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.1059'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.04%'
# Row 14
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001301'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '2.85%'
# Row 15
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005898'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-1.91%'
# Row 16
$ws.Range("B16").Value = 'HotbitToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.004412'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.24%'
# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '1.65%'
# Row 18
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.416'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.30%'
# Row 19
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3344'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.19%'
# Row 20
$ws.Range("B20").Value = 'MCDex'
$ws.Range("C20").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.571'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-16.11%'
# Row 21
$ws.Range("B21").Value = 'ProBitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1371'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-1.37%'
# Row 22
$ws.Range("B22").Value = 'ZBToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2577'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '0.41%'
# Row 23
$ws.Range("B23").Value = 'CoinExToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04159'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '1.61%'
# Row 24
$ws.Range("B24").Value = 'BitKan'
$ws.Range("C24").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001234'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-5.36%'
# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '2.02%'
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02685'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '0.16%'
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05560'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-1.74%'
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007819'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '2.78%'
# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-0.81%'
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.007337'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-2.44%'
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002113'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '8.18%'
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007832'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-5.36%'
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006833'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-2.92%'
# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.41%'
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003492'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-1.18%'
# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.41%'
# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.41%'
